$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'22.448.28"
$ws.Cells.Item(2, 5).Value = "  +0.29%  "

$ws.Cells.Item(3, 4).Value = "'1.568.64"
$ws.Cells.Item(3, 5).Value = "  +0.06%  "

$ws.Cells.Item(4, 4).Value = "'1.003"
$ws.Cells.Item(4, 5).Value = "  -0.02%  "

$ws.Cells.Item(5, 4).Value = "'1.003"
$ws.Cells.Item(5, 5).Value = "  -0.10%  "

$ws.Cells.Item(6, 4).Value = "'290.26"
$ws.Cells.Item(6, 5).Value = "  -0.19%  "

$ws.Cells.Item(7, 4).Value = "'0.3690"
$ws.Cells.Item(7, 5).Value = "  -1.40%  "

$ws.Cells.Item(8, 4).Value = "'49.78"
$ws.Cells.Item(8, 5).Value = "  +1.52%  "

$ws.Cells.Item(9, 4).Value = "'0.3367"
$ws.Cells.Item(9, 5).Value = "  -0.57%  "

$ws.Cells.Item(10, 4).Value = "'1.145"
$ws.Cells.Item(10, 5).Value = "  +1.27%  "

$ws.Cells.Item(11, 4).Value = "'0.07532"
$ws.Cells.Item(11, 5).Value = "  -0.04%  "

$ws.Cells.Item(12, 5).Value = "  -0.08%  "

$ws.Cells.Item(13, 4).Value = "'21.10"
$ws.Cells.Item(13, 5).Value = "  +1.01%  "

$ws.Cells.Item(14, 5).Value = "  +1.46%  "

$ws.Cells.Item(15, 4).Value = "'6.964"
$ws.Cells.Item(15, 5).Value = "  +1.19%  "

$ws.Cells.Item(16, 4).Value = "'1.571.88"
$ws.Cells.Item(16, 5).Value = "  +0.46%  "

$ws.Cells.Item(17, 4).Value = "'0.00001119"
$ws.Cells.Item(17, 5).Value = "  -0.28%  "

$ws.Cells.Item(18, 4).Value = "'90.36"
$ws.Cells.Item(18, 5).Value = "  +0.73%  "

$ws.Cells.Item(19, 4).Value = "'0.06768"
$ws.Cells.Item(19, 5).Value = "  +0.52%  "

$ws.Cells.Item(20, 4).Value = "'1.003"
$ws.Cells.Item(20, 5).Value = "  -0.08%  "

$ws.Cells.Item(21, 4).Value = "'6.359"
$ws.Cells.Item(21, 5).Value = "  +3.03%  "

$ws.Cells.Item(22, 4).Value = "'16.36"
$ws.Cells.Item(22, 5).Value = "  -0.46%  "

$ws.Cells.Item(23, 4).Value = "'12.21"
$ws.Cells.Item(23, 5).Value = "  +2.60%  "

$ws.Cells.Item(24, 4).Value = "'22.455.32"
$ws.Cells.Item(24, 5).Value = "  +0.38%  "

$ws.Cells.Item(25, 4).Value = "'2.380"
$ws.Cells.Item(25, 5).Value = "  +0.03%  "

$ws.Cells.Item(26, 4).Value = "'2.647"
$ws.Cells.Item(26, 5).Value = "  -2.24%  "

$ws.Cells.Item(27, 4).Value = "'19.98"
$ws.Cells.Item(27, 5).Value = "  -0.27%  "

$ws.Cells.Item(28, 4).Value = "'148.97"
$ws.Cells.Item(28, 5).Value = "  +1.05%  "

$ws.Cells.Item(29, 4).Value = "'5.061"
$ws.Cells.Item(29, 5).Value = "  +0.57%  "

$ws.Cells.Item(30, 4).Value = "'124.86"
$ws.Cells.Item(30, 5).Value = "  -0.32%  "

$ws.Cells.Item(31, 4).Value = "'1.750.53"
$ws.Cells.Item(31, 5).Value = "  +0.54%  "

$ws.Cells.Item(32, 4).Value = "'1.060"
$ws.Cells.Item(32, 5).Value = "  +7.67%  "

$ws.Cells.Item(33, 4).Value = "'6.177"
$ws.Cells.Item(33, 5).Value = "  +2.79%  "

$ws.Cells.Item(34, 4).Value = "'2.011"
$ws.Cells.Item(34, 5).Value = "  -0.40%  "

$ws.Cells.Item(35, 4).Value = "'9.783"
$ws.Cells.Item(35, 5).Value = "  -2.30%  "

$ws.Cells.Item(36, 4).Value = "'0.08332"
$ws.Cells.Item(36, 5).Value = "  -1.59%  "

$ws.Cells.Item(37, 2).Value = "TrustWalletToken"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(37, 4).Value = "'1.357"
$ws.Cells.Item(37, 5).Value = "  -4.35%  "

$ws.Cells.Item(38, 2).Value = "VeChain"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(38, 4).Value = "'0.02466"
$ws.Cells.Item(38, 5).Value = "  -0.42%  "

$ws.Cells.Item(39, 4).Value = "'0.2297"
$ws.Cells.Item(39, 5).Value = "  +0.92%  "

$ws.Cells.Item(40, 4).Value = "'0.06554"
$ws.Cells.Item(40, 5).Value = "  +1.87%  "

$ws.Cells.Item(41, 4).Value = "'5.415"
$ws.Cells.Item(41, 5).Value = "  +0.61%  "

$ws.Cells.Item(42, 4).Value = "'11.19"
$ws.Cells.Item(42, 5).Value = "  +1.22%  "

$ws.Cells.Item(43, 5).Value = "  -0.74%  "

$ws.Cells.Item(44, 4).Value = "'14.11"
$ws.Cells.Item(44, 5).Value = "  +1.36%  "

$ws.Cells.Item(45, 5).Value = "  -0.08%  "

$ws.Cells.Item(46, 4).Value = "'3.805"
$ws.Cells.Item(46, 5).Value = "  +0.11%  "

$ws.Cells.Item(47, 4).Value = "'0.5843"
$ws.Cells.Item(47, 5).Value = "  -0.59%  "

$ws.Cells.Item(48, 5).Value = "  +0.47%  "

$ws.Cells.Item(49, 4).Value = "'127.75"
$ws.Cells.Item(49, 5).Value = "  +2.63%  "

$ws.Cells.Item(50, 5).Value = "  -1.60%  "

$ws.Cells.Item(51, 4).Value = "'0.07300"
$ws.Cells.Item(51, 5).Value = "  -0.29%  "
